$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G3').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G4').Value = 'mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range('H4').Value = '103/216'
$ws.Range('G5').Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G7').Value = 'aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range('G8').Value = 'dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range('G13').Value = 'ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G16').Value = 'aya.hanafy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G18').Value = 'nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G19').Value = 'asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G21').Value = 'Amera.a.saad@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G23').Value = 'aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range('G24').Value = 'dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range('G29').Value = 'ManaratAleslam-Zefan@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G32').Value = 'aya.hanafy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G34').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G35').Value = 'gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, heba@med.asu.edu.eg'
$ws.Range('G50').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G51').Value = 'gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, heba@med.asu.edu.eg'
$ws.Range('G66').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G67').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G69').Value = 'menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Range('G77').Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range('G80').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G82').Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G83').Value = 'asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G85').Value = 'menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg'
$ws.Range('G93').Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range('G96').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G98').Value = 'nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G103').Value = 'aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range('G114').Value = 'nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G116').Value = 'mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range('G119').Value = 'aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'

# Percentage-looking text values (e.g. "55.9%") get auto-converted to numeric
# percent values by Excel if assigned directly via .Value. Use a scratch cell
# pre-formatted as Text, then PasteSpecial (values only) into the target cell so
# the target's original style/number format (General, inlineStr-equivalent text)
# is preserved while the literal text is carried over untouched.

$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '55.9%'
$ws.Range('ZZ1').Copy()
$ws.Range('L10').PasteSpecial(-4163)
$ws.Range('ZZ1').Clear()

$ws.Range('ZZ1').NumberFormat = '@'
$ws.Range('ZZ1').Value = '57.6%'
$ws.Range('ZZ1').Copy()
$ws.Range('S15').PasteSpecial(-4163)
$ws.Range('ZZ1').Clear()
